# Add a copy of the "LKT 8HED3" sheet, named "LKT 8HED3A", at the end of
# the workbook. This mirrors (in Excel) duplicating the sheet tab via
# right-click > Move or Copy... > Create a copy, placing it after the
# last existing sheet, then renaming the new tab.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("LKT 8HED3")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Copy places the new sheet immediately after $lastSheet and makes it active.
$source.Copy($null, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "LKT 8HED3A"
